$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new column before column C.
#    Old columns C (Statistical test) and D (Notes) shift to D and E.
# ------------------------------------------------------------------
$ws.Columns("C").Insert()

# ------------------------------------------------------------------
# 2. Retarget the old "C2:C11 not-blank" conditional format onto the
#    shifted D2:D11 range and fix up its formulas to reference D2
#    instead of C2. Do this before any other conditional-format
#    changes so the range queries below aren't confused by overlap.
# ------------------------------------------------------------------
$oldBlankRules = $ws.Range("C2:C11").FormatConditions
$blankNotBlank = $oldBlankRules.Item(1)
$blankIsBlank  = $oldBlankRules.Item(2)
$blankNotBlank.Formula1 = "=LEN(TRIM(D2))>0"
$blankIsBlank.Formula1  = "=LEN(TRIM(D2))=0"
$blankNotBlank.ModifyAppliesToRange($ws.Range("D2:D11"))

# ------------------------------------------------------------------
# 3. Extend the existing B3:B11 "Yes/No" rule to also cover column C.
# ------------------------------------------------------------------
$bRules = $ws.Range("B3:B11").FormatConditions
$bRules.Item(1).ModifyAppliesToRange($ws.Range("B3:C11"))

# ------------------------------------------------------------------
# 4. Add a brand-new Yes/No rule for the lone cell C2.
# ------------------------------------------------------------------
$cRules = $ws.Range("C2").FormatConditions
$cEqual = $cRules.Add(1, 3, '"Yes"')
$cEqual.Font.Color = 0x006100
$cEqual.Interior.Color = 0xCEEFC6

$cNotEqual = $cRules.Add(1, 4, '"Yes"')
$cNotEqual.Font.Color = 0x06009C
$cNotEqual.Interior.Color = 0xCEC7FF

# ------------------------------------------------------------------
# 5. Fill in the new column C header + "No" data for every data row,
#    and clear its inherited formatting so it keeps the default style.
# ------------------------------------------------------------------
$ws.Range("C1").Value = "Complete (N recorded for generating data)?"

$ws.Range("C2:C11").ClearFormats()
$ws.Range("C2:C11").Value = "No"

# ------------------------------------------------------------------
# 6. Update column B ("Complete (small error and notes)?") so every
#    row now reads "Yes".
# ------------------------------------------------------------------
$ws.Range("B2").Value = "Yes"
$ws.Range("B6").Value = "Yes"
$ws.Range("B7").Value = "Yes"

# ------------------------------------------------------------------
# 7. Update the notes column (now E) with new note text for rows 2 and 7.
#    (E6 already holds the untouched rich-text note that shifted with
#    the column insert.)
# ------------------------------------------------------------------
$ws.Range("E2").Value = "Several cases where it's not really group polarization: either polarity crosses neutral point or initial opinion is at neutral point; both of these seem to contradict the assumption of bias conservation."
$ws.Range("E7").Value = "All results are BS? Paper uses ANOVA to compare shifts, but no significance testing of shifts"

# ------------------------------------------------------------------
# 8. Update the sheet view: scroll to show column C first and select E11.
# ------------------------------------------------------------------
$ws.Range("E11").Select()
$excel.ActiveWindow.ScrollColumn = 3
